$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.687.51"
$ws.Range("E2").Value = "  +0.31%  "

$ws.Range("D3").Value = "3.647.50"
$ws.Range("E3").Value = "  +1.94%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.66"
$ws.Range("E5").Value = "  +0.80%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.85"
$ws.Range("E6").Value = "  +19.30%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "655.69"
$ws.Range("E7").Value = "  -0.18%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.422"
$ws.Range("E8").Value = "  +4.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.07"
$ws.Range("E9").Value = "  +1.93%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.00"
$ws.Range("E10").Value = "  -0.02%  "

$ws.Range("D11").Value = "3.643.01"
$ws.Range("E11").Value = "  +1.83%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "44.39"
$ws.Range("E12").Value = "  +2.46%  "

$ws.Range("E13").Value = "  +1.05%  "

$ws.Range("E14").Value = "  +2.23%  "

$ws.Range("D15").Value = "4.326.83"
$ws.Range("E15").Value = "  +1.98%  "

$ws.Range("D16").Value = "96.552.39"
$ws.Range("E16").Value = "  +0.27%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000259"
$ws.Range("E17").Value = "  +0.41%  "

$ws.Range("D18").Value = "3.647.51"
$ws.Range("E18").Value = "  +2.21%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.95"
$ws.Range("E19").Value = "  +2.72%  "

$ws.Range("E20").Value = "  +0.05%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.38"
$ws.Range("E21").Value = "  +3.43%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.535"
$ws.Range("E22").Value = "  +8.75%  "

$ws.Range("E23").Value = "  +0.42%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "512.95"
$ws.Range("E24").Value = "  +0.23%  "

$ws.Range("E25").Value = "  +2.87%  "

$ws.Range("E26").Value = "  +1.21%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "101.16"
$ws.Range("E27").Value = "  +4.91%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "13.11"
$ws.Range("E28").Value = "  +2.48%  "

$ws.Range("E29").Value = "  +12.06%  "

$ws.Range("E30").Value = "  +1.75%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.89"
$ws.Range("E31").Value = "  +3.68%  "

$ws.Range("E32").Value = "  -0.14%  "

$ws.Range("E33").Value = "  +1.23%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "33.17"
$ws.Range("E34").Value = "  +5.00%  "

$ws.Range("E35").Value = "  +0.16%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.76"
$ws.Range("E36").Value = "  +8.72%  "

$ws.Range("E37").Value = "  +3.82%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.84"
$ws.Range("E38").Value = "  +3.14%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "616.49"
$ws.Range("E39").Value = "  +1.29%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "42.97"
$ws.Range("E40").Value = "  +25.97%  "

$ws.Range("E41").Value = "  +4.72%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.954"
$ws.Range("E42").Value = "  +5.36%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.94"
$ws.Range("E43").Value = "  +6.65%  "

$ws.Range("E44").Value = "  +0.00%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.14"
$ws.Range("E45").Value = "  +7.20%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0442"
$ws.Range("E46").Value = "  +5.95%  "

$ws.Range("E47").Value = "  +1.75%  "

$ws.Range("E48").Value = "  +0.28%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.412"
$ws.Range("E49").Value = "  +17.97%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.62"
$ws.Range("E50").Value = "  +5.18%  "

$ws.Range("E51").Value = "  +1.54%  "
